$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 9).Value = "sv"
$ws.Cells.Item(2, 10).Value = "Statement-opinion"
$ws.Cells.Item(8, 9).Value = "b"
$ws.Cells.Item(8, 10).Value = "Acknowledge (Backchannel)"
$ws.Cells.Item(9, 9).Value = "aa"
$ws.Cells.Item(9, 10).Value = "Agree/Accept"
$ws.Cells.Item(13, 9).Value = "sv"
$ws.Cells.Item(13, 10).Value = "Statement-opinion"
$ws.Cells.Item(15, 9).Value = "sv"
$ws.Cells.Item(15, 10).Value = "Statement-opinion"
$ws.Cells.Item(19, 9).Value = "sd"
$ws.Cells.Item(19, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(32, 9).Value = "sv"
$ws.Cells.Item(32, 10).Value = "Statement-opinion"
$ws.Cells.Item(34, 9).Value = "sv"
$ws.Cells.Item(34, 10).Value = "Statement-opinion"
$ws.Cells.Item(36, 9).Value = "sd"
$ws.Cells.Item(36, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(37, 9).Value = "sv"
$ws.Cells.Item(37, 10).Value = "Statement-opinion"
$ws.Cells.Item(47, 9).Value = "aa"
$ws.Cells.Item(47, 10).Value = "Agree/Accept"
$ws.Cells.Item(48, 9).Value = "%"
$ws.Cells.Item(48, 10).Value = "Uninterpretable"
$ws.Cells.Item(49, 9).Value = "ba"
$ws.Cells.Item(49, 10).Value = "Appreciation"
$ws.Cells.Item(50, 9).Value = "b"
$ws.Cells.Item(50, 10).Value = "Acknowledge (Backchannel)"
$ws.Cells.Item(67, 9).Value = "aa"
$ws.Cells.Item(67, 10).Value = "Agree/Accept"
$ws.Cells.Item(71, 9).Value = "sd"
$ws.Cells.Item(71, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(72, 9).Value = "aa"
$ws.Cells.Item(72, 10).Value = "Agree/Accept"
$ws.Cells.Item(76, 9).Value = "aa"
$ws.Cells.Item(76, 10).Value = "Agree/Accept"
$ws.Cells.Item(77, 9).Value = "sv"
$ws.Cells.Item(77, 10).Value = "Statement-opinion"
$ws.Cells.Item(78, 9).Value = "aa"
$ws.Cells.Item(78, 10).Value = "Agree/Accept"
$ws.Cells.Item(113, 9).Value = "b"
$ws.Cells.Item(113, 10).Value = "Acknowledge (Backchannel)"
$ws.Cells.Item(119, 9).Value = "sd"
$ws.Cells.Item(119, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(125, 9).Value = "aa"
$ws.Cells.Item(125, 10).Value = "Agree/Accept"
$ws.Cells.Item(132, 9).Value = "sv"
$ws.Cells.Item(132, 10).Value = "Statement-opinion"
$ws.Cells.Item(134, 9).Value = "aa"
$ws.Cells.Item(134, 10).Value = "Agree/Accept"
$ws.Cells.Item(142, 9).Value = "aa"
$ws.Cells.Item(142, 10).Value = "Agree/Accept"
$ws.Cells.Item(143, 9).Value = "ba"
$ws.Cells.Item(143, 10).Value = "Appreciation"
$ws.Cells.Item(144, 9).Value = "sd"
$ws.Cells.Item(144, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(158, 9).Value = "aa"
$ws.Cells.Item(158, 10).Value = "Agree/Accept"
$ws.Cells.Item(172, 9).Value = "aa"
$ws.Cells.Item(172, 10).Value = "Agree/Accept"
$ws.Cells.Item(173, 9).Value = "aa"
$ws.Cells.Item(173, 10).Value = "Agree/Accept"
$ws.Cells.Item(177, 9).Value = "sv"
$ws.Cells.Item(177, 10).Value = "Statement-opinion"
$ws.Cells.Item(184, 9).Value = "%"
$ws.Cells.Item(184, 10).Value = "Uninterpretable"
$ws.Cells.Item(190, 9).Value = "sd"
$ws.Cells.Item(190, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(195, 9).Value = "sv"
$ws.Cells.Item(195, 10).Value = "Statement-opinion"
$ws.Cells.Item(196, 9).Value = "sv"
$ws.Cells.Item(196, 10).Value = "Statement-opinion"
$ws.Cells.Item(199, 9).Value = "sv"
$ws.Cells.Item(199, 10).Value = "Statement-opinion"
$ws.Cells.Item(209, 9).Value = "sd"
$ws.Cells.Item(209, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(211, 9).Value = "aa"
$ws.Cells.Item(211, 10).Value = "Agree/Accept"
$ws.Cells.Item(215, 9).Value = "sd"
$ws.Cells.Item(215, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(216, 9).Value = "sv"
$ws.Cells.Item(216, 10).Value = "Statement-opinion"
$ws.Cells.Item(221, 9).Value = "sd"
$ws.Cells.Item(221, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(227, 9).Value = "sv"
$ws.Cells.Item(227, 10).Value = "Statement-opinion"
$ws.Cells.Item(231, 9).Value = "aa"
$ws.Cells.Item(231, 10).Value = "Agree/Accept"
$ws.Cells.Item(232, 9).Value = "aa"
$ws.Cells.Item(232, 10).Value = "Agree/Accept"
$ws.Cells.Item(240, 9).Value = "aa"
$ws.Cells.Item(240, 10).Value = "Agree/Accept"
$ws.Cells.Item(244, 9).Value = "sd"
$ws.Cells.Item(244, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(248, 9).Value = "sd"
$ws.Cells.Item(248, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(249, 9).Value = "sd"
$ws.Cells.Item(249, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(256, 9).Value = "sv"
$ws.Cells.Item(256, 10).Value = "Statement-opinion"
$ws.Cells.Item(260, 9).Value = "aa"
$ws.Cells.Item(260, 10).Value = "Agree/Accept"
$ws.Cells.Item(265, 9).Value = "ba"
$ws.Cells.Item(265, 10).Value = "Appreciation"
$ws.Cells.Item(277, 9).Value = "sd"
$ws.Cells.Item(277, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(281, 9).Value = "sd"
$ws.Cells.Item(281, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(282, 9).Value = "sd"
$ws.Cells.Item(282, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(284, 9).Value = "%"
$ws.Cells.Item(284, 10).Value = "Uninterpretable"
$ws.Cells.Item(305, 9).Value = "sd"
$ws.Cells.Item(305, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(306, 9).Value = "sd"
$ws.Cells.Item(306, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(309, 9).Value = "sd"
$ws.Cells.Item(309, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(322, 9).Value = "sd"
$ws.Cells.Item(322, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(325, 9).Value = "sv"
$ws.Cells.Item(325, 10).Value = "Statement-opinion"
$ws.Cells.Item(331, 9).Value = "b"
$ws.Cells.Item(331, 10).Value = "Acknowledge (Backchannel)"
$ws.Cells.Item(355, 9).Value = "b"
$ws.Cells.Item(355, 10).Value = "Acknowledge (Backchannel)"
